$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 275
$ws1.Range("F5").Value = 280
$ws1.Range("F6").Value = 1083
$ws1.Range("F7").Value = 1423
$ws1.Range("F8").Value = 585
$ws1.Range("F9").Value = 103
$ws1.Range("F10").Value = 745
$ws1.Range("F12").Value = 147
$ws1.Range("F13").Value = 122
$ws1.Range("F15").Value = 1337
$ws1.Range("F16").Value = 103
$ws1.Range("F17").Value = 95
$ws1.Range("F18").Value = 274
$ws1.Range("F19").Value = 5213
$ws1.Range("F21").Value = 32
$ws1.Range("F22").Value = 209
$ws1.Range("F23").Value = 18
$ws1.Range("F24").Value = 5745
$ws1.Range("F25").Value = 58
$ws1.Range("F29").Value = 14325
$ws1.Range("F30").Value = 1425
$ws1.Range("F32").Value = 97
$ws1.Range("F34").Value = 2934
$ws1.Range("F35").Value = 598
$ws1.Range("F36").Value = 4185
$ws1.Range("F37").Value = 129
$ws1.Range("F38").Value = 355
$ws1.Range("F39").Value = 110

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 275
$ws4.Range("F5").Value = 280
$ws4.Range("F6").Value = 1083
$ws4.Range("F7").Value = 1423
$ws4.Range("F8").Value = 585
$ws4.Range("F9").Value = 103
$ws4.Range("F10").Value = 745
$ws4.Range("F12").Value = 147
$ws4.Range("F13").Value = 122
$ws4.Range("F15").Value = 1337
$ws4.Range("F16").Value = 103
$ws4.Range("F17").Value = 95
$ws4.Range("F18").Value = 274
$ws4.Range("F20").Value = 5213
$ws4.Range("F23").Value = 32
$ws4.Range("F24").Value = 209
$ws4.Range("F25").Value = 18
$ws4.Range("F27").Value = 5745
$ws4.Range("F28").Value = 58
$ws4.Range("F32").Value = 14325
$ws4.Range("F33").Value = 1425
$ws4.Range("F35").Value = 97
$ws4.Range("F37").Value = 2934
$ws4.Range("F38").Value = 598
$ws4.Range("F39").Value = 4185
$ws4.Range("F40").Value = 129
$ws4.Range("F41").Value = 355
$ws4.Range("F42").Value = 110
